# The deck originally carries the "Integral" theme (ppt/theme/theme1.xml,
# used by the one-and-only slide master) while the notes master points at
# a spare, unused "Office Theme" theme part (ppt/theme/theme2.xml).
#
# The authored change swaps the two themes' colour schemes around: the
# slides start using the stock "Office Theme" palette, while the
# previously-active "Integral" palette is pushed onto the (otherwise
# inert) second theme part.
#
# Apply that by rewriting the 12 theme colours on the presentation's
# live theme (Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink,
# FollowedHyperlink) from "Integral" to the standard Office palette.

function ConvertTo-BgrColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in ThemeColorScheme index order:
# 1 Dark1, 2 Light1, 3 Dark2, 4 Light2, 5-10 Accent1-6, 11 Hyperlink,
# 12 FollowedHyperlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-BgrColor $officeThemeColors[$i - 1]
}
